# The edit shuffles the data rows (2-41) of the sheet: for each row r, the
# values in columns D, L, M, N, O, P, R, S are replaced by the values that
# were originally in another row (the mapping below), while columns
# A, B, C, E, F, G, H, I, J, K, Q, T stay as-is (they are constant across
# all rows anyway).
#
# Because this is a permutation (not a simple shift), we must snapshot all
# the "before" values first, then write the "after" values from the
# snapshot - otherwise we would overwrite source data before it is read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> source row (new row r gets old row's data from $map[r])
$map = @{
    2 = 7;  3 = 5;  4 = 31; 5 = 40; 6 = 11; 7 = 35; 8 = 12; 9 = 13; 10 = 32;
    11 = 24; 12 = 38; 13 = 30; 14 = 18; 15 = 34; 16 = 21; 17 = 28; 18 = 33;
    19 = 36; 20 = 39; 21 = 26; 22 = 22; 23 = 16; 24 = 8; 25 = 10; 26 = 3;
    27 = 17; 28 = 20; 29 = 2; 30 = 27; 31 = 14; 32 = 29; 33 = 41; 34 = 6;
    35 = 37; 36 = 4;  37 = 23; 38 = 9;  39 = 19; 40 = 25; 41 = 15
}

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the "before" values for every relevant column/row.
$snapshot = @{}
foreach ($col in $cols) {
    $snapshot[$col] = @{}
    for ($r = 2; $r -le 41; $r++) {
        $snapshot[$col][$r] = $ws.Range("$col$r").Value()
    }
}

# Write the "after" values using the snapshot so the permutation is applied
# consistently, independent of write order.
foreach ($col in $cols) {
    for ($r = 2; $r -le 41; $r++) {
        $src = $map[$r]
        $ws.Range("$col$r").Value = $snapshot[$col][$src]
    }
}
